$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original "2026/12/29" entry (row 736) is being replaced by a new
# "2026/02/01" reading; all the following rows (old 736:777) shift down by
# one to make room, ending up at 737:778.
$ws.Rows.Item(736).Insert()

# Populate the newly inserted row 736 with the new reading.
# Date-looking text must be forced to Text so Excel doesn't reinterpret
# "2026/02/01" as a date serial number; ClearFormats() afterwards drops the
# temporary Text number format so the cell ends up with the default style,
# matching every other date cell in the column.
$dateCell = $ws.Cells.Item(736, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/01"
$dateCell.ClearFormats()

$ws.Cells.Item(736, 2).Value = "日"
$ws.Cells.Item(736, 3).Value = 16
$ws.Cells.Item(736, 4).Value = 20
